$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "57.667.16"
$ws.Range("E2").Value = "  -0.08%  "

# Row 3
$ws.Range("D3").Value = "3.065.23"
$ws.Range("E3").Value = "  +1.44%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.04%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "516.16"
$ws.Range("E5").Value = "  +1.08%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.66"
$ws.Range("E6").Value = "  +0.20%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.04%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.436"
$ws.Range("E8").Value = "  -0.32%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.23"
$ws.Range("E9").Value = "  -4.87%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.110"
$ws.Range("E10").Value = "  -0.98%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.379"
$ws.Range("E11").Value = "  +2.96%  "

# Row 12
$ws.Range("D12").Value = "3.584.75"
$ws.Range("E12").Value = "  +1.48%  "

# Row 13
$ws.Range("E13").Value = "  -3.26%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.92"
$ws.Range("E14").Value = "  +0.89%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000169"
$ws.Range("E15").Value = "  +2.82%  "

# Row 16
$ws.Range("D16").Value = "57.738.57"
$ws.Range("E16").Value = "  +0.16%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.21"
$ws.Range("E17").Value = "  -1.24%  "

# Row 18
$ws.Range("D18").Value = "3.064.73"
$ws.Range("E18").Value = "  +1.52%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.51"
$ws.Range("E19").Value = "  +4.33%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.20"
$ws.Range("E20").Value = "  +2.43%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "332.10"
$ws.Range("E21").Value = "  +0.60%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("E22").Value = "  +0.08%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.509"
$ws.Range("E23").Value = "  +1.50%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.30"
$ws.Range("E24").Value = "  +0.78%  "

# Row 25
$ws.Range("D25").Value = "3.185.88"
$ws.Range("E25").Value = "  +1.50%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  +0.05%  "

# Row 27
$ws.Range("E27").Value = "  -2.60%  "

# Row 28
$ws.Range("D28").Value = "0.0₃0910"
$ws.Range("E28").Value = "  -1.57%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.78"
$ws.Range("E29").Value = "  -1.03%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.27"
$ws.Range("E30").Value = "  -1.62%  "

# Row 31
$ws.Range("E31").Value = "  -0.21%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.22"
$ws.Range("E32").Value = "  +1.15%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.91"
$ws.Range("E33").Value = "  +0.94%  "

# Row 34
$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.69"
$ws.Range("E34").Value = "  -1.82%  "

# Row 35
$ws.Range("B35").Value = "Monero"
$ws.Range("C35").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "153.92"
$ws.Range("E35").Value = "  -0.83%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.93"
$ws.Range("E36").Value = "  +0.38%  "

# Row 37
$ws.Range("E37").Value = "  -0.62%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "25.07"
$ws.Range("E38").Value = "  +1.76%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0683"
$ws.Range("E39").Value = "  +0.49%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.25"
$ws.Range("E40").Value = "  -1.05%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.91"
$ws.Range("E41").Value = "  +0.87%  "

# Row 42
$ws.Range("B42").Value = "Mantle"
$ws.Range("C42").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.670"
$ws.Range("E42").Value = "  +2.89%  "

# Row 43
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.999"
$ws.Range("E43").Value = "  +0.02%  "

# Row 44
$ws.Range("E44").Value = "  -1.00%  "

# Row 45
$ws.Range("D45").Value = "2.208.95"
$ws.Range("E45").Value = "  -1.12%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.12"
$ws.Range("E46").Value = "  +1.12%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.955"
$ws.Range("E47").Value = "  -3.48%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0244"
$ws.Range("E48").Value = "  +1.69%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "20.30"
$ws.Range("E49").Value = "  +3.70%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0175"
$ws.Range("E50").Value = "  +7.55%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.185"
$ws.Range("E51").Value = "  -0.26%  "
